$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update time_taken (column F) timestamps on the "data" sheet ---
$ws.Range("F2").Value  = "2021-10-05 14:21:39.917656"
$ws.Range("F3").Value  = "2021-10-05 14:21:39.917664"
$ws.Range("F4").Value  = "2021-10-05 14:21:39.917667"
$ws.Range("F5").Value  = "2021-10-05 14:21:39.917670"
$ws.Range("F6").Value  = "2021-10-05 14:21:39.917673"
$ws.Range("F7").Value  = "2021-10-05 14:21:39.917676"
$ws.Range("F8").Value  = "2021-10-05 14:21:39.917679"
$ws.Range("F9").Value  = "2021-10-05 14:21:39.917681"
$ws.Range("F10").Value = "2021-10-05 14:21:39.917684"
$ws.Range("F11").Value = "2021-10-05 14:21:39.917687"
$ws.Range("F12").Value = "2021-10-05 14:21:39.917690"

# --- Add a new "metadata" worksheet positioned after "data" ---
$meta = $wb.Worksheets.Add($null, $ws)
$meta.Name = "metadata"

# Reuse the bold/centered header style already used on the "data" sheet
# (column B1:F1) instead of creating brand new style entries.
$ws.Range("B1:F1").Copy()
$meta.Range("B1:F1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("F1").Copy()
$meta.Range("G1").PasteSpecial(-4122)      # xlPasteFormats
$ws.Range("A2").Copy()
$meta.Range("A2").PasteSpecial(-4122)      # xlPasteFormats

# Header row
$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

# Data row
$meta.Range("A2").Value = 0
$meta.Range("B2").Value = "Multiple Epiphyseal Dysplasia"
$meta.Range("C2").Value = 211
$meta.Range("D2").NumberFormat = "@"
$meta.Range("D2").Value = "1.2"
$meta.Range("E2").Value = "2017-11-05T02:37:20.234212Z"
$meta.Range("F2").Value = "2021-10-05 14:21:39.914242"
$meta.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/211/?format=json"
